$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 787
$ws.Range("J28").Value = 1324.3334
$ws.Range("L28").Value = 1324.3334
$ws.Range("N28").Value = -2294.3334
$ws.Range("H62").Value = 3820.1052
$ws.Range("I62").Value = 3132.5
$ws.Range("J62").Value = 4998.857
$ws.Range("K62").Value = 3132.5
$ws.Range("L62").Value = 4998.857
$ws.Range("M62").Value = -2508.5
$ws.Range("N62").Value = -6246.857
$ws.Range("H64").Value = 3754.5454
$ws.Range("I64").Value = 3450
$ws.Range("J64").Value = 3928.5715
$ws.Range("K64").Value = 3450
$ws.Range("L64").Value = 3928.5715
$ws.Range("M64").Value = -3202
$ws.Range("N64").Value = -4424.5715
$ws.Range("H65").Value = 3820.1052
$ws.Range("I65").Value = 3132.5
$ws.Range("J65").Value = 4998.857
$ws.Range("K65").Value = 15662.5
$ws.Range("L65").Value = 24994.285
$ws.Range("M65").Value = -12542.5
$ws.Range("N65").Value = -31234.285
$ws.Range("H67").Value = 3754.5454
$ws.Range("I67").Value = 3450
$ws.Range("J67").Value = 3928.5715
$ws.Range("K67").Value = 3450
$ws.Range("L67").Value = 3928.5715
$ws.Range("M67").Value = -2592
$ws.Range("N67").Value = -5644.5715
$ws.Range("H80").Value = 6768885.5
$ws.Range("I80").Value = 617.17645
$ws.Range("K80").Value = 1851.52935
$ws.Range("M80").Value = -853.5293500000002
$ws.Range("H83").Value = 6768885.5
$ws.Range("I83").Value = 617.17645
$ws.Range("K83").Value = 5554.58805
$ws.Range("M83").Value = -562.5880500000003
$ws.Range("H92").Value = 71429380
$ws.Range("I92").Value = 100000820
$ws.Range("K92").Value = 100000820
$ws.Range("M92").Value = -99999572
$ws.Range("H115").Value = 741.3333
$ws.Range("I115").Value = 741.3333
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2223.9999
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -656.9998999999998
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 932
$ws.Range("I118").Value = 642.5
$ws.Range("J118").Value = 2090
$ws.Range("K118").Value = 1927.5
$ws.Range("L118").Value = 6270
$ws.Range("M118").Value = -270.5
$ws.Range("N118").Value = -9584
$ws.Range("H137").Value = 1793.8064
$ws.Range("I137").Value = 1467.5
$ws.Range("J137").Value = 2245.6155
$ws.Range("K137").Value = 4402.5
$ws.Range("L137").Value = 6736.8465
$ws.Range("M137").Value = -1852.5
$ws.Range("N137").Value = -11836.8465
$ws.Range("H138").Value = 2250.6057
$ws.Range("I138").Value = 2251.4614
$ws.Range("J138").Value = 2250.4138
$ws.Range("K138").Value = 6754.3842
$ws.Range("L138").Value = 6751.241399999999
$ws.Range("M138").Value = -1614.3842
$ws.Range("N138").Value = -17031.2414
$ws.Range("H141").Value = 1060.5686
$ws.Range("I141").Value = 825.1111
$ws.Range("J141").Value = 2826.5
$ws.Range("K141").Value = 2475.3333
$ws.Range("L141").Value = 8479.5
$ws.Range("M141").Value = 2704.6667
$ws.Range("N141").Value = -18839.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 163.33333
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H32").Value = 5144.4907
$ws.Range("I32").Value = 5089.577
$ws.Range("K32").Value = 5089.577
$ws.Range("M32").Value = -4802.577
$ws.Range("H61").Value = 1194.2307
$ws.Range("I61").Value = 1074.5834
$ws.Range("K61").Value = 1074.5834
$ws.Range("M61").Value = -862.5834
$ws.Range("H74").Value = 83334330
$ws.Range("I74").Value = 90910070
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 90910070
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = -90909196
$ws.Range("N74").Value = -2948
$ws.Range("H77").Value = 83334330
$ws.Range("I77").Value = 90910070
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 454550350
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = -454545982
$ws.Range("N77").Value = -14736
$ws.Range("H97").Value = 1374.7084
$ws.Range("I97").Value = 1458.7778
$ws.Range("J97").Value = 1122.5
$ws.Range("K97").Value = 1458.7778
$ws.Range("L97").Value = 1122.5
$ws.Range("M97").Value = -962.7778000000001
$ws.Range("N97").Value = -2114.5
$ws.Range("H110").Value = 749
$ws.Range("I110").Value = 749
$ws.Range("K110").Value = 749
$ws.Range("M110").Value = 1296
$ws.Range("H122").Value = 1766.1316
$ws.Range("I122").Value = 1418.1428
$ws.Range("J122").Value = 2740.5
$ws.Range("K122").Value = 4254.428400000001
$ws.Range("L122").Value = 8221.5
$ws.Range("M122").Value = -1804.428400000001
$ws.Range("N122").Value = -13121.5
$ws.Range("H132").Value = 26343.951
$ws.Range("I132").Value = 1501.36
$ws.Range("K132").Value = 4504.08
$ws.Range("M132").Value = -1974.08
$ws.Range("H136").Value = 1194.2307
$ws.Range("I136").Value = 1074.5834
$ws.Range("K136").Value = 3223.7502
$ws.Range("M136").Value = -673.7501999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H59").Value = 45000
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -46694
$ws.Range("H86").Value = 1410.9736
$ws.Range("I86").Value = 1296.8077
$ws.Range("J86").Value = 1658.3334
$ws.Range("K86").Value = 1296.8077
$ws.Range("L86").Value = 1658.3334
$ws.Range("M86").Value = -173.8077000000001
$ws.Range("N86").Value = -3904.3334
$ws.Range("H89").Value = 1410.9736
$ws.Range("I89").Value = 1296.8077
$ws.Range("J89").Value = 1658.3334
$ws.Range("K89").Value = 6484.038500000001
$ws.Range("L89").Value = 8291.666999999999
$ws.Range("M89").Value = -868.0385000000006
$ws.Range("N89").Value = -19523.667
$ws.Range("H94").Value = 444.73685
$ws.Range("I94").Value = 458.33334
$ws.Range("J94").Value = 200
$ws.Range("K94").Value = 458.33334
$ws.Range("L94").Value = 200
$ws.Range("M94").Value = -7.333340000000021
$ws.Range("N94").Value = -1102
$ws.Range("H99").Value = 1491.3334
$ws.Range("I99").Value = 1374.4445
$ws.Range("K99").Value = 1374.4445
$ws.Range("M99").Value = 123.5554999999999
$ws.Range("H134").Value = 5124.68
$ws.Range("I134").Value = 5345.591
$ws.Range("K134").Value = 16036.773
$ws.Range("M134").Value = -13501.773

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 217.83333
$ws.Range("I7").Value = 221.4
$ws.Range("K7").Value = 221.4
$ws.Range("M7").Value = -108.4
$ws.Range("H16").Value = 1011.5
$ws.Range("I16").Value = 752
$ws.Range("J16").Value = 1271
$ws.Range("K16").Value = 752
$ws.Range("L16").Value = 1271
$ws.Range("M16").Value = -465
$ws.Range("N16").Value = -1845
$ws.Range("H31").Value = 15589.4
$ws.Range("I31").Value = 19296.863
$ws.Range("J31").Value = 5393.875
$ws.Range("K31").Value = 19296.863
$ws.Range("L31").Value = 5393.875
$ws.Range("M31").Value = -19001.863
$ws.Range("N31").Value = -5983.875
$ws.Range("H34").Value = 15589.4
$ws.Range("I34").Value = 19296.863
$ws.Range("J34").Value = 5393.875
$ws.Range("K34").Value = 19296.863
$ws.Range("L34").Value = 5393.875
$ws.Range("M34").Value = -19094.863
$ws.Range("N34").Value = -5797.875
$ws.Range("H69").Value = 4792.4287
$ws.Range("I69").Value = 2931.8
$ws.Range("K69").Value = 2931.8
$ws.Range("M69").Value = -2182.8
$ws.Range("H72").Value = 4792.4287
$ws.Range("I72").Value = 2931.8
$ws.Range("K72").Value = 8795.400000000001
$ws.Range("M72").Value = -5051.400000000001
$ws.Range("H99").Value = 17245352
$ws.Range("I99").Value = 3735.45
$ws.Range("K99").Value = 3735.45
$ws.Range("M99").Value = -2237.45
$ws.Range("H107").Value = 1923.7273
$ws.Range("I107").Value = 2750
$ws.Range("J107").Value = 1740.1111
$ws.Range("K107").Value = 2750
$ws.Range("L107").Value = 1740.1111
$ws.Range("M107").Value = -830
$ws.Range("N107").Value = -5580.1111
$ws.Range("H113").Value = 1011.5
$ws.Range("I113").Value = 752
$ws.Range("J113").Value = 1271
$ws.Range("K113").Value = 752
$ws.Range("L113").Value = 1271
$ws.Range("M113").Value = 1418
$ws.Range("N113").Value = -5611
$ws.Range("H126").Value = 17245352
$ws.Range("I126").Value = 3735.45
$ws.Range("K126").Value = 11206.35
$ws.Range("M126").Value = -8736.349999999999
$ws.Range("H132").Value = 10889.741
$ws.Range("I132").Value = 12774.441
$ws.Range("J132").Value = 3522.2727
$ws.Range("K132").Value = 38323.323
$ws.Range("L132").Value = 10566.8181
$ws.Range("M132").Value = -35793.323
$ws.Range("N132").Value = -15626.8181
$ws.Range("H134").Value = 739.5
$ws.Range("I134").Value = 604.5217
$ws.Range("J134").Value = 1084.4445
$ws.Range("K134").Value = 1813.5651
$ws.Range("L134").Value = 3253.3335
$ws.Range("M134").Value = 721.4349
$ws.Range("N134").Value = -8323.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3075.75
$ws.Range("I3").Value = 1518.7
$ws.Range("K3").Value = 4556.1
$ws.Range("M3").Value = -4444.1
$ws.Range("H9").Value = 1001
$ws.Range("J9").Value = 1001
$ws.Range("L9").Value = 3003
$ws.Range("N9").Value = -3451
$ws.Range("H92").Value = 12500374
$ws.Range("I92").Value = 25000356
$ws.Range("J92").Value = 390.6
$ws.Range("K92").Value = 75001068
$ws.Range("L92").Value = 1171.8
$ws.Range("M92").Value = -74999820
$ws.Range("N92").Value = -3667.8
$ws.Range("H122").Value = 834.8570999999999
$ws.Range("I122").Value = 376.25
$ws.Range("J122").Value = 942.7646999999999
$ws.Range("K122").Value = 3386.25
$ws.Range("L122").Value = 8484.882299999999
$ws.Range("M122").Value = -936.25
$ws.Range("N122").Value = -13384.8823
$ws.Range("H131").Value = 170320.77
$ws.Range("I131").Value = 553.3333
$ws.Range("J131").Value = 179415.45
$ws.Range("K131").Value = 1659.9999
$ws.Range("L131").Value = 538246.3500000001
$ws.Range("M131").Value = 3380.0001
$ws.Range("N131").Value = -548326.3500000001
$ws.Range("H133").Value = 4333.1665
$ws.Range("I133").Value = 2000
$ws.Range("J133").Value = 5999.7144
$ws.Range("K133").Value = 6000
$ws.Range("L133").Value = 17999.1432
$ws.Range("M133").Value = -940
$ws.Range("N133").Value = -28119.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3729
$ws.Range("I80").Value = 3350
$ws.Range("J80").Value = 3880.6
$ws.Range("K80").Value = 3350
$ws.Range("L80").Value = 3880.6
$ws.Range("M80").Value = -2352
$ws.Range("N80").Value = -5876.6
$ws.Range("H83").Value = 3729
$ws.Range("I83").Value = 3350
$ws.Range("J83").Value = 3880.6
$ws.Range("K83").Value = 16750
$ws.Range("L83").Value = 19403
$ws.Range("M83").Value = -11758
$ws.Range("N83").Value = -29387
$ws.Range("H97").Value = 996.8125
$ws.Range("I97").Value = 1011.46155
$ws.Range("J97").Value = 933.3333
$ws.Range("K97").Value = 1011.46155
$ws.Range("L97").Value = 933.3333
$ws.Range("M97").Value = -515.46155
$ws.Range("N97").Value = -1925.3333
$ws.Range("H113").Value = 2314.6365
$ws.Range("I113").Value = 1694
$ws.Range("J113").Value = 3059.4
$ws.Range("K113").Value = 1694
$ws.Range("L113").Value = 3059.4
$ws.Range("M113").Value = 476
$ws.Range("N113").Value = -7399.4
$ws.Range("H123").Value = 10325.833
$ws.Range("J123").Value = 10325.833
$ws.Range("L123").Value = 10325.833
$ws.Range("N123").Value = -15225.833
$ws.Range("H132").Value = 16807.416
$ws.Range("I132").Value = 2947.037
$ws.Range("K132").Value = 8841.110999999999
$ws.Range("M132").Value = -6311.110999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5846.091
$ws.Range("I61").Value = 2586.7144
$ws.Range("J61").Value = 11550
$ws.Range("K61").Value = 2586.7144
$ws.Range("L61").Value = 11550
$ws.Range("M61").Value = -2384.7144
$ws.Range("N61").Value = -11954
$ws.Range("H68").Value = 2409.9
$ws.Range("I68").Value = 2450
$ws.Range("J68").Value = 2249.5
$ws.Range("K68").Value = 2450
$ws.Range("L68").Value = 2249.5
$ws.Range("M68").Value = -1701
$ws.Range("N68").Value = -3747.5
$ws.Range("H71").Value = 2409.9
$ws.Range("I71").Value = 2450
$ws.Range("J71").Value = 2249.5
$ws.Range("K71").Value = 12250
$ws.Range("L71").Value = 11247.5
$ws.Range("M71").Value = -8506
$ws.Range("N71").Value = -18735.5
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 5846.091
$ws.Range("I113").Value = 2586.7144
$ws.Range("J113").Value = 11550
$ws.Range("K113").Value = 2586.7144
$ws.Range("L113").Value = 11550
$ws.Range("M113").Value = -416.7143999999998
$ws.Range("N113").Value = -15890
$ws.Range("H132").Value = 2573.7368
$ws.Range("I132").Value = 2221.6428
$ws.Range("J132").Value = 3559.6
$ws.Range("K132").Value = 6664.928400000001
$ws.Range("L132").Value = 10678.8
$ws.Range("M132").Value = -4134.928400000001
$ws.Range("N132").Value = -15738.8
$ws.Range("H136").Value = 1384.8276
$ws.Range("I136").Value = 1130.5555
$ws.Range("J136").Value = 1800.909
$ws.Range("K136").Value = 3391.6665
$ws.Range("L136").Value = 5402.727000000001
$ws.Range("M136").Value = -841.6664999999998
$ws.Range("N136").Value = -10502.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3399.3809
$ws.Range("I62").Value = 3057.25
$ws.Range("J62").Value = 3855.5557
$ws.Range("K62").Value = 3057.25
$ws.Range("L62").Value = 3855.5557
$ws.Range("M62").Value = -2433.25
$ws.Range("N62").Value = -5103.5557
$ws.Range("H65").Value = 3399.3809
$ws.Range("I65").Value = 3057.25
$ws.Range("J65").Value = 3855.5557
$ws.Range("K65").Value = 15286.25
$ws.Range("L65").Value = 19277.7785
$ws.Range("M65").Value = -12166.25
$ws.Range("N65").Value = -25517.7785
$ws.Range("H113").Value = 6761757
$ws.Range("I113").Value = 6666.6665
$ws.Range("J113").Value = 27027028
$ws.Range("K113").Value = 19999.9995
$ws.Range("L113").Value = 81081084
$ws.Range("M113").Value = -17829.9995
$ws.Range("N113").Value = -81085424
$ws.Range("H126").Value = 1888.9333
$ws.Range("I126").Value = 1175.4546
$ws.Range("J126").Value = 3851
$ws.Range("K126").Value = 3526.3638
$ws.Range("L126").Value = 11553
$ws.Range("M126").Value = -1056.3638
$ws.Range("N126").Value = -16493
$ws.Range("H132").Value = 813.7105
$ws.Range("I132").Value = 581.19446
$ws.Range("K132").Value = 1743.58338
$ws.Range("M132").Value = 786.41662
